# "error solve ifrs list"
# Recompute the IFRS financial figures on the company_list sheet: rows 2-6
# get corrected (much smaller) numbers in columns D:AJ, a handful of cells
# in those rows that no longer apply are cleared, and rows 7-9 (the
# 2019/12(E)-2021/12(E) forecast years) lose all of their D:AJ figures,
# leaving only the A/B/C label columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 188
$ws.Range("E2").Value = -60
$ws.Range("F2").Value = -60
$ws.Range("G2").Value = -187
$ws.Range("H2").Value = -187
$ws.Range("I2").Value = -187
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 488
$ws.Range("L2").Value = 110
$ws.Range("M2").Value = 378
$ws.Range("N2").Value = 378
$ws.Range("O2").ClearContents()
$ws.Range("P2").Value = 268
$ws.Range("Q2").Value = -31
$ws.Range("R2").Value = -344
$ws.Range("S2").Value = 343
$ws.Range("T2").Value = 36
$ws.Range("U2").Value = -66
$ws.Range("V2").Value = 83
$ws.Range("W2").Value = -31.88
$ws.Range("X2").Value = -99.45
$ws.Range("Y2").Value = -66.73
$ws.Range("Z2").Value = -45.5
$ws.Range("AA2").Value = 29.15
$ws.Range("AB2").Value = 47.3
$ws.Range("AC2").Value = -449
$ws.Range("AD2").Value = -1.61
$ws.Range("AE2").Value = 469
$ws.Range("AF2").Value = 1.54
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 81040022

# Row 3
$ws.Range("D3").Value = 92
$ws.Range("E3").Value = -15
$ws.Range("F3").Value = -15
$ws.Range("G3").Value = 15
$ws.Range("H3").Value = 12
$ws.Range("I3").Value = 12
$ws.Range("J3").ClearContents()
$ws.Range("K3").Value = 480
$ws.Range("L3").Value = 30
$ws.Range("M3").Value = 450
$ws.Range("N3").Value = 450
$ws.Range("O3").ClearContents()
$ws.Range("P3").Value = 291
$ws.Range("Q3").Value = 2
$ws.Range("R3").Value = 292
$ws.Range("S3").Value = -35
$ws.Range("T3").Value = 6
$ws.Range("U3").Value = -4
$ws.Range("V3").ClearContents()
$ws.Range("W3").Value = -15.95
$ws.Range("X3").Value = 12.86
$ws.Range("Y3").Value = 2.85
$ws.Range("Z3").Value = 2.43
$ws.Range("AA3").Value = 6.61
$ws.Range("AB3").Value = 58.81
$ws.Range("AC3").Value = 14
$ws.Range("AD3").Value = 121.34
$ws.Range("AE3").Value = 514
$ws.Range("AF3").Value = 3.37
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 88057219

# Row 4
$ws.Range("D4").Value = 98
$ws.Range("E4").Value = -16
$ws.Range("F4").Value = -16
$ws.Range("G4").Value = -23
$ws.Range("H4").Value = -23
$ws.Range("I4").Value = -23
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 521
$ws.Range("L4").Value = 29
$ws.Range("M4").Value = 492
$ws.Range("N4").Value = 491
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 318
$ws.Range("Q4").Value = -35
$ws.Range("R4").Value = -249
$ws.Range("S4").Value = 73
$ws.Range("T4").Value = 2
$ws.Range("U4").Value = -37
$ws.Range("V4").ClearContents()
$ws.Range("W4").Value = -16.73
$ws.Range("X4").Value = -23.85
$ws.Range("Y4").Value = -4.96
$ws.Range("Z4").Value = -4.68
$ws.Range("AA4").Value = 5.85
$ws.Range("AB4").Value = 60.12
$ws.Range("AC4").Value = -26
$ws.Range("AD4").Value = -34.11
$ws.Range("AE4").Value = 515
$ws.Range("AF4").Value = 1.75
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 96279524

# Row 5
$ws.Range("D5").Value = 77
$ws.Range("E5").Value = -29
$ws.Range("F5").Value = -29
$ws.Range("G5").Value = -27
$ws.Range("H5").Value = -27
$ws.Range("I5").Value = -27
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 688
$ws.Range("L5").Value = 31
$ws.Range("M5").Value = 657
$ws.Range("N5").Value = 656
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 518
$ws.Range("Q5").Value = -13
$ws.Range("R5").Value = -29
$ws.Range("S5").Value = 195
$ws.Range("T5").Value = 27
$ws.Range("U5").Value = -39
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = -38.11
$ws.Range("X5").Value = -34.99
$ws.Range("Y5").Value = -4.64
$ws.Range("Z5").Value = -4.46
$ws.Range("AA5").Value = 4.75
$ws.Range("AB5").Value = 30.81
$ws.Range("AC5").Value = -26
$ws.Range("AD5").Value = -16.82
$ws.Range("AE5").Value = 510
$ws.Range("AF5").Value = 0.85
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 129484149

# Row 6
$ws.Range("D6").Value = 64
$ws.Range("E6").Value = -14
$ws.Range("F6").Value = -14
$ws.Range("G6").Value = -108
$ws.Range("H6").Value = -108
$ws.Range("I6").Value = -108
$ws.Range("K6").Value = 576
$ws.Range("L6").Value = 24
$ws.Range("M6").Value = 551
$ws.Range("N6").Value = 551
$ws.Range("P6").Value = 518
$ws.Range("Q6").Value = -30
$ws.Range("R6").Value = -5
$ws.Range("S6").Value = -3
$ws.Range("T6").Value = 2
$ws.Range("U6").Value = -32
$ws.Range("V6").Value = 0
$ws.Range("W6").Value = -22.09
$ws.Range("X6").Value = -168.25
$ws.Range("Y6").Value = -17.96
$ws.Range("Z6").Value = -17.16
$ws.Range("AA6").Value = 4.42
$ws.Range("AB6").Value = 9.01
$ws.Range("AC6").Value = -84
$ws.Range("AD6").Value = -3.92
$ws.Range("AE6").Value = 428
$ws.Range("AF6").Value = 0.77
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 129484149

# Row 7: clear all data columns, keep A/B/C
$ws.Range("D7:AJ7").ClearContents()

# Row 8: clear all data columns, keep A/B/C
$ws.Range("D8:AJ8").ClearContents()

# Row 9: clear all data columns, keep A/B/C
$ws.Range("D9:AJ9").ClearContents()
